{"js": "// Update the date line and the 25 two-digit-by-two-digit multiplication\n// prompts throughout the document (table cells). Each old value is unique\n// in the document, so a simple search-and-replace per pair is safe and\n// order-independent.\nconst replacements = [\n  [\"2024-10-09 Wednesday\", \"2024-10-10 Thursday\"],\n  [\"39\u00d746=\", \"74\u00d732=\"],\n  [\"76\u00d736=\", \"75\u00d719=\"],\n  [\"85\u00d736=\", \"91\u00d721=\"],\n  [\"93\u00d722=\", \"91\u00d796=\"],\n  [\"72\u00d791=\", \"79\u00d797=\"],\n  [\"71\u00d785=\", \"52\u00d724=\"],\n  [\"35\u00d717=\", \"96\u00d785=\"],\n  [\"95\u00d783=\", \"58\u00d796=\"],\n  [\"60\u00d724=\", \"66\u00d768=\"],\n  [\"38\u00d720=\", \"33\u00d799=\"],\n  [\"59\u00d758=\", \"72\u00d792=\"],\n  [\"93\u00d791=\", \"48\u00d792=\"],\n  [\"86\u00d780=\", \"42\u00d781=\"],\n  [\"94\u00d720=\", \"79\u00d790=\"],\n  [\"59\u00d744=\", \"38\u00d714=\"],\n  [\"62\u00d759=\", \"56\u00d739=\"],\n  [\"24\u00d780=\", \"80\u00d771=\"],\n  [\"70\u00d711=\", \"75\u00d757=\"],\n  [\"39\u00d770=\", \"43\u00d732=\"],\n  [\"82\u00d750=\", \"87\u00d715=\"],\n  [\"25\u00d715=\", \"20\u00d799=\"],\n  [\"36\u00d737=\", \"62\u00d744=\"],\n  [\"91\u00d745=\", \"34\u00d791=\"],\n  [\"13\u00d726=\", \"45\u00d711=\"],\n  [\"73\u00d727=\", \"16\u00d752=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the 25 two-digit-by-two-digit multiplication\n# prompts throughout the document (table cells). Each old value is unique\n# in the document, so a simple Find/Replace per pair is safe.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-10-09 Wednesday\", \"2024-10-10 Thursday\"),\n    @(\"39\u00d746=\", \"74\u00d732=\"),\n    @(\"76\u00d736=\", \"75\u00d719=\"),\n    @(\"85\u00d736=\", \"91\u00d721=\"),\n    @(\"93\u00d722=\", \"91\u00d796=\"),\n    @(\"72\u00d791=\", \"79\u00d797=\"),\n    @(\"71\u00d785=\", \"52\u00d724=\"),\n    @(\"35\u00d717=\", \"96\u00d785=\"),\n    @(\"95\u00d783=\", \"58\u00d796=\"),\n    @(\"60\u00d724=\", \"66\u00d768=\"),\n    @(\"38\u00d720=\", \"33\u00d799=\"),\n    @(\"59\u00d758=\", \"72\u00d792=\"),\n    @(\"93\u00d791=\", \"48\u00d792=\"),\n    @(\"86\u00d780=\", \"42\u00d781=\"),\n    @(\"94\u00d720=\", \"79\u00d790=\"),\n    @(\"59\u00d744=\", \"38\u00d714=\"),\n    @(\"62\u00d759=\", \"56\u00d739=\"),\n    @(\"24\u00d780=\", \"80\u00d771=\"),\n    @(\"70\u00d711=\", \"75\u00d757=\"),\n    @(\"39\u00d770=\", \"43\u00d732=\"),\n    @(\"82\u00d750=\", \"87\u00d715=\"),\n    @(\"25\u00d715=\", \"20\u00d799=\"),\n    @(\"36\u00d737=\", \"62\u00d744=\"),\n    @(\"91\u00d745=\", \"34\u00d791=\"),\n    @(\"13\u00d726=\", \"45\u00d711=\"),\n    @(\"73\u00d727=\", \"16\u00d752=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $found = $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
